# Updates cryptos list values (price / volume(1h)) to reflect the latest
# scrape, as produced by the scheduled GitHub Actions job.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper that writes a literal text value into a cell without letting Excel's
# automatic "number/date detection" mangle numeric-looking strings (e.g.
# "1.00" silently becoming the number 1, or "492.64" becoming a noisy
# floating point value). We briefly force the cell to Text format, assign
# the value, then restore the cell style so no stray formatting is left
# behind on the cell.
function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$Val
    )
    $rng = $ws.Range($CellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $Val
    $rng.Style = "Normal"
}

# --- Price (column D) updates ---
Set-TextValue "D2"  "53.793.91"
Set-TextValue "D3"  "2.250.73"
Set-TextValue "D5"  "492.64"
Set-TextValue "D7"  "0.999"
Set-TextValue "D8"  "0.522"
Set-TextValue "D9"  "0.0942"
Set-TextValue "D12" "4.72"
Set-TextValue "D13" "2.649.12"
Set-TextValue "D14" "22.46"
Set-TextValue "D15" "53.775.89"
Set-TextValue "D17" "2.250.07"
Set-TextValue "D18" "10.15"
Set-TextValue "D19" "4.11"
Set-TextValue "D20" "301.20"
Set-TextValue "D21" "6.26"
Set-TextValue "D22" "1.00"
Set-TextValue "D23" "60.58"
Set-TextValue "D24" "0.999"
Set-TextValue "D26" "7.23"
Set-TextValue "D27" "170.86"
Set-TextValue "D31" "1.06"
Set-TextValue "D33" "17.69"
Set-TextValue "D34" "0.999"
Set-TextValue "D35" "0.933"
Set-TextValue "D36" "1.18"
Set-TextValue "D37" "3.67"
Set-TextValue "D38" "0.369"
Set-TextValue "D39" "1.38"
Set-TextValue "D40" "3.33"
Set-TextValue "D41" "123.93"
Set-TextValue "D42" "4.75"
Set-TextValue "D43" "0.0486"
Set-TextValue "D45" "0.538"
Set-TextValue "D46" "237.18"
Set-TextValue "D47" "0.368"
Set-TextValue "D49" "10.74"
Set-TextValue "D50" "16.06"
Set-TextValue "D51" "4.61"

# --- Volume(1h) (column E) updates ---
Set-TextValue "E2"  "  -1.72%  "
Set-TextValue "E3"  "  -2.07%  "
Set-TextValue "E4"  "  +0.17%  "
Set-TextValue "E5"  "  -0.94%  "
Set-TextValue "E6"  "  -0.75%  "
Set-TextValue "E7"  "  +0.05%  "
Set-TextValue "E8"  "  -1.43%  "
Set-TextValue "E9"  "  -0.90%  "
Set-TextValue "E10" "  +0.52%  "
Set-TextValue "E11" "  +2.75%  "
Set-TextValue "E12" "  +1.49%  "
Set-TextValue "E13" "  -2.11%  "
Set-TextValue "E14" "  +2.89%  "
Set-TextValue "E15" "  -1.37%  "
Set-TextValue "E16" "  -0.77%  "
Set-TextValue "E17" "  -2.10%  "
Set-TextValue "E18" "  +0.88%  "
Set-TextValue "E19" "  -0.24%  "
Set-TextValue "E20" "  -0.58%  "
Set-TextValue "E21" "  -3.42%  "
Set-TextValue "E22" "  +0.16%  "
Set-TextValue "E23" "  -3.47%  "
Set-TextValue "E24" "  -0.15%  "
Set-TextValue "E25" "  -2.87%  "
Set-TextValue "E26" "  +1.56%  "
Set-TextValue "E27" "  +0.70%  "
Set-TextValue "E28" "  -0.95%  "
Set-TextValue "E31" "  -1.51%  "
Set-TextValue "E32" "  -0.01%  "
Set-TextValue "E33" "  +0.13%  "
Set-TextValue "E34" "  +0.05%  "
Set-TextValue "E35" "  +6.96%  "
Set-TextValue "E36" "  -1.49%  "
Set-TextValue "E37" "  -0.82%  "
Set-TextValue "E38" "  -2.03%  "
Set-TextValue "E39" "  -2.33%  "
Set-TextValue "E40" "  -0.81%  "
Set-TextValue "E41" "  -3.21%  "
Set-TextValue "E42" "  -2.09%  "
Set-TextValue "E43" "  +0.01%  "
Set-TextValue "E44" "  -0.66%  "
Set-TextValue "E45" "  -1.85%  "
Set-TextValue "E46" "  -2.10%  "
Set-TextValue "E47" "  -1.84%  "
Set-TextValue "E48" "  -0.66%  "
Set-TextValue "E49" "  +0.25%  "
Set-TextValue "E50" "  -2.36%  "
Set-TextValue "E51" "  -1.02%  "

# --- Rows 29/30 swapped places in the ranking (PEPE <-> Aptos) ---
Set-TextValue "B29" "Aptos"
Set-TextValue "C29" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D29" "5.88"
Set-TextValue "E29" "  -0.34%  "

Set-TextValue "B30" "PEPE"
Set-TextValue "C30" "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue "D30" "0.0₃0680"
Set-TextValue "E30" "  -2.21%  "
